# Rename the two duplicate-named embedded logo pictures that live in the
# document's headers/footers:
#   - the Pearson Edexcel logo (currently "image2.png")  -> "image1.png"
#   - the BTEC orange logo     (currently "image1.jpg")  -> "image2.jpg"
#
# The pictures are inline drawings inside the first-page header/footer and
# the default footer, so we walk every section's Headers/Footers collections
# (covering wdHeaderFooterPrimary/FirstPage/EvenPages) and rename each
# InlineShape we find by matching on its (stable) AlternativeText, rather
# than relying on brittle positional indexing.

$d = $word.ActiveDocument

function Rename-LogoShapes($range) {
    if ($range -eq $null) { return }
    $shapes = $range.InlineShapes
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        $desc = $shp.AlternativeText

        if ($desc -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shp.Name = "image1.png"
        } elseif ($desc -eq "BTec_Logo-Orange") {
            $shp.Name = "image2.jpg"
        }
    }
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($i = 1; $i -le $sec.Headers.Count; $i++) {
        $h = $sec.Headers.Item($i)
        if ($h.Exists) {
            Rename-LogoShapes $h.Range
        }
    }

    for ($i = 1; $i -le $sec.Footers.Count; $i++) {
        $f = $sec.Footers.Item($i)
        if ($f.Exists) {
            Rename-LogoShapes $f.Range
        }
    }
}

# Also cover any matching pictures that might live in the document body.
Rename-LogoShapes $d.Content
